$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.078.41"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "'2.118.78"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'346.60"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.5193"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "'0.4476"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").Value = "'54.25"
$ws.Range("E9").Value = "  +3.46%  "

$ws.Range("D10").Value = "'0.09364"
$ws.Range("E10").Value = "  -0.77%  "

$ws.Range("D11").Value = "'1.184"
$ws.Range("E11").Value = "  +0.50%  "

$ws.Range("D12").Value = "'25.43"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "'8.691"
$ws.Range("E13").Value = "  +7.30%  "

$ws.Range("D14").Value = "'6.979"
$ws.Range("E14").Value = "  +3.55%  "

$ws.Range("D15").Value = "'2.110.59"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").Value = "'102.52"
$ws.Range("E16").Value = "  +3.01%  "

$ws.Range("D17").Value = "'0.00001171"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").Value = "'21.61"
$ws.Range("E19").Value = "  +4.43%  "

$ws.Range("D20").Value = "'0.06701"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").Value = "'6.303"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "'30.111.65"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").Value = "'12.75"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").Value = "'2.333"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "'2.360.81"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'22.17"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.542"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'162.71"
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'134.24"
$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.159"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.794"
$ws.Range("E32").Value = "  +10.15%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.1057"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'6.283"
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'6.698"
$ws.Range("E35").Value = "  +8.36%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'3.970"
$ws.Range("E36").Value = "  +0.57%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'10.74"
$ws.Range("E37").Value = "  +5.86%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02640"
$ws.Range("E38").Value = "  +2.44%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06870"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.7143"
$ws.Range("E40").Value = "  +2.66%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'12.73"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2254"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.328"
$ws.Range("E43").Value = "  +2.09%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6915"
$ws.Range("E44").Value = "  +3.55%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'14.69"
$ws.Range("E45").Value = "  +3.35%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.395"
$ws.Range("E46").Value = "  +4.66%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.006"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").Value = "'3.633"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'1.261"
$ws.Range("E49").Value = "  +7.64%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000356"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.229"
$ws.Range("E51").Value = "  +0.50%  "
